# Auto-generated Excel COM-interop script
# Commit: "Updated symbol list on Wed Feb 15 16:43:18 UTC 2023 with GitHub Actions"
#
# The sheet is a scraped crypto-price table (Coin / Link / Price / Volume(1h) /
# Data / Hora). This refresh run:
#   - updates Price (D) and Volume(1h) (E) text for many existing rows,
#   - inserts a new row (GateToken) at row 6 and shifts the rows that used to
#     follow it down by one (each row keeps its own old Coin/Link, just moved
#     down one slot) through row 17,
#   - drops the last row (CoinExToken/BitKan/... tail) off that shifted block.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'301.85"
$ws.Range("E2").Value = "'1.98%"

# Row 3
$ws.Range("D3").Value = "'43.79"
$ws.Range("E3").Value = "'5.98%"

# Row 4
$ws.Range("D4").Value = "'5.081"
$ws.Range("E4").Value = "'0.83%"

# Row 5
$ws.Range("D5").Value = "'0.07696"
$ws.Range("E5").Value = "'3.44%"

# Row 6
$ws.Range("B6").Value = 'GateToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D6").Value = "'4.422"
$ws.Range("E6").Value = "'1.49%"

# Row 7
$ws.Range("B7").Value = 'FTXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D7").Value = "'1.620"
$ws.Range("E7").Value = "'3.15%"

# Row 8
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").Value = "'1.047"
$ws.Range("E8").Value = "'13.46%"

# Row 9
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").Value = "'0.1274"
$ws.Range("E9").Value = "'7.97%"

# Row 10
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = "'0.1884"
$ws.Range("E10").Value = "'3.29%"

# Row 11
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = "'0.09164"
$ws.Range("E11").Value = "'4.19%"

# Row 12
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = "'0.04167"
$ws.Range("E12").Value = "'-2.98%"

# Row 13
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = "'0.1049"
$ws.Range("E13").Value = "'-0.21%"

# Row 14
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = "'0.001273"
$ws.Range("E14").Value = "'-0.63%"

# Row 15
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").Value = "'0.005768"
$ws.Range("E15").Value = "'-3.32%"

# Row 16
$ws.Range("B16").Value = 'UpBots'
$ws.Range("C16").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("D16").Value = "'0.007430"
$ws.Range("E16").Value = "'1,895.26%"

# Row 17
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = "'3.343"
$ws.Range("E17").Value = "'-0.39%"

# Row 18
$ws.Range("E18").Value = "'-2.68%"

# Row 19
$ws.Range("D19").Value = "'0.3349"
$ws.Range("E19").Value = "'1.40%"

# Row 20
$ws.Range("D20").Value = "'8.654"
$ws.Range("E20").Value = "'10.01%"

# Row 21
$ws.Range("D21").Value = "'0.1400"
$ws.Range("E21").Value = "'2.11%"

# Row 22
$ws.Range("D22").Value = "'0.3180"
$ws.Range("E22").Value = "'7.14%"

# Row 23
$ws.Range("D23").Value = "'0.04176"
$ws.Range("E23").Value = "'3.62%"

# Row 24
$ws.Range("D24").Value = "'0.001286"
$ws.Range("E24").Value = "'1.32%"

# Row 25
$ws.Range("D25").Value = "'0.004425"
$ws.Range("E25").Value = "'14.85%"

# Row 26
$ws.Range("D26").Value = "'0.0001348"
$ws.Range("E26").Value = "'9.56%"

# Row 38
$ws.Range("E38").Value = "'3.83%"

# Row 39
$ws.Range("D39").Value = "'0.05301"
$ws.Range("E39").Value = "'1.98%"

# Row 40
$ws.Range("D40").Value = "'0.005940"
$ws.Range("E40").Value = "'-11.46%"

# Row 41
$ws.Range("D41").Value = "'0.007657"
$ws.Range("E41").Value = "'-1.57%"

# Row 42
$ws.Range("D42").Value = "'0.1347"
$ws.Range("E42").Value = "'2.42%"

# Row 43
$ws.Range("D43").Value = "'0.007382"
$ws.Range("E43").Value = "'-0.01%"

# Row 44
$ws.Range("D44").Value = "'0.007546"
$ws.Range("E44").Value = "'-3.36%"

# Row 45
$ws.Range("D45").Value = "'0.3004"
$ws.Range("E45").Value = "'-6.56%"

# Row 46
$ws.Range("D46").Value = "'0.00006661"
$ws.Range("E46").Value = "'5.76%"

# Row 47
$ws.Range("D47").Value = "'0.00000000749"
$ws.Range("E47").Value = "'-0.17%"

# Row 48
$ws.Range("D48").Value = "'0.04077"
$ws.Range("E48").Value = "'-11.47%"

# Row 49
$ws.Range("E49").Value = "'-0.04%"

# Row 50
$ws.Range("D50").Value = "'0.00002098"
$ws.Range("E50").Value = "'-0.17%"

# Row 51
$ws.Range("D51").Value = "'0.0001998"
$ws.Range("E51").Value = "'-0.17%"
